$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row: number correct/right mark value
$ws.Range("B11").Value = 5

# Update "Total" row: total score and the "correct/total" summary text
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"
